# Add an explicit "do not start this paragraph on a new page"
# (pageBreakBefore = False) paragraph setting to every paragraph in the
# document body, and to the built-in heading/title paragraph styles.

$d = $word.ActiveDocument

# 1) Every paragraph currently in the document gets an explicit
#    <w:pageBreakBefore w:val="0"/> in its pPr.
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# 2) The built-in heading / title paragraph styles get the same explicit
#    setting recorded on their paragraph format.
$styleNames = @(
    "Heading 1",
    "Heading 2",
    "Heading 3",
    "Heading 4",
    "Heading 5",
    "Heading 6",
    "Title",
    "Subtitle"
)

foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    $s.ParagraphFormat.PageBreakBefore = 0
}
